# Fixed a bug with IsAlwaysGen in chgSymbols
# The data rows (A2:F25) are reshuffled to new row positions; each row keeps
# its own set of 6 values (symbol id + reel1..reel5 weights), only the row
# position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(101,9,30,15,60,15),
    @(701,3,90,45,97,15),
    @(601,9,60,67,60,42),
    @(801,3,67,65,52,45),
    @(902,1,0,0,0,0),
    @(501,9,52,30,75,45),
    @(201,9,30,15,45,30),
    @(301,6,45,30,60,45),
    @(401,9,48,67,75,45),
    @(1201,2,10,10,10,10),
    @(901,16,15,45,60,60),
    @(1001,18,30,75,60,72),
    @(1202,2,10,10,10,10),
    @(1203,3,15,15,15,15),
    @(1101,0,15,30,30,0),
    @(2,0,2,2,2,2),
    @(502,0,4,0,0,0),
    @(802,0,4,5,4,0),
    @(1,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(402,0,0,4,0,0),
    @(602,0,0,4,0,9),
    @(702,0,0,0,4,0),
    @(1002,0,0,0,0,9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
